# GithubActions-WorkflowsInYourRepo.pptx - "Fixed code changes in slides"
#
# 1) Slide 10 ("How GitHub Actions Work"): trim the last bullet and drop the
#    trailing empty paragraph.
# 2) Slides 21, 22, 23 ("Creating a CI/CD Pipeline"): rename the default
#    branch in the sample YAML from `master` to `main`.
# 3) Slide 27 ("Creating a CI/CD Pipeline"): drop the hyphen in the sample
#    app-name value.

$p = $ppt.ActivePresentation

# --- Slide 10: "Each step of the job is executed on the same runner to allow" ---
$slide10 = $p.Slides.Item(10)
$shape10 = $slide10.Shapes.Item(2)
$tr10 = $shape10.TextFrame.TextRange
$para10_3 = $tr10.Paragraphs(3, 1)
$para10_3.Runs(1, 1).Text = "Each step of the job is executed on the same runner"
$para10_4 = $tr10.Paragraphs(4, 1)
$para10_4.Delete()

# --- Slides 21/22/23: "    branches: [ master ]" -> "    branches: [ main ]" ---
foreach ($slideIdx in 21, 22, 23) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(2)
    $tr = $shape.TextFrame.TextRange
    for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
        $para = $tr.Paragraphs($i, 1)
        if ($para.Text.TrimEnd("`r") -eq "    branches: [ master ]") {
            $para.Runs(1, 1).Text = "    branches: [ main ]"
            break
        }
    }
}

# --- Slide 27: "app-name: 'GitHubActions-PCC2021'" -> "app-name: 'GitHubActionsPCC2021'" ---
$leftSingleQuote = [char]0x2018
$slide27 = $p.Slides.Item(27)
$shape27 = $slide27.Shapes.Item(2)
$tr27 = $shape27.TextFrame.TextRange
for ($i = 1; $i -le $tr27.Paragraphs().Count; $i++) {
    $para = $tr27.Paragraphs($i, 1)
    if ($para.Text -like "*app-name:*GitHubActions-PCC2021*") {
        $para.Runs(1, 1).Text = "        app-name: " + $leftSingleQuote + "GitHubActionsPCC2021'"
        break
    }
}
